$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge runs: "Виконання та захист " + "практич" + "них робіт"
#    -> single run "Виконання та захист практичних робіт"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Виконання та захист практичних робіт", $true, $false, $false, $false, $false, $true, 1, $false, "Виконання та захист практичних робіт", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge runs: "«" + "2" + "» х " + "7" + " = " + "14"
#    -> single run "«2» х 7 = 14"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("«2» х 7 = 14", $true, $false, $false, $false, $false, $true, 1, $false, "«2» х 7 = 14", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Merge runs: "«" + "4" + "» х" + "7" + " = " + "28"
#    -> single run "«4» х7 = 28"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("«4» х7 = 28", $true, $false, $false, $false, $false, $true, 1, $false, "«4» х7 = 28", 2) | Out-Null

# ---------------------------------------------------------------------------
# Helper: split a single run of text into 4 runs around two 1-char "digit"
# substitutions while keeping formatting identical (a transient Bold
# toggle forces the run boundary to persist instead of being re-merged).
#   searchText   - the current, single-run text to locate
#   digitPos     - 0-based offset of the single-character token to replace
#                  (e.g. the "6" in «6» ...)
#   newDigit     - replacement text for that single character
#   tailStart    - 0-based offset where the trailing numeric token starts
#                  (e.g. the "18" in ... = 18)
#   tailEnd      - 0-based offset just past the trailing numeric token
#   newTail      - replacement text for the trailing numeric token
# ---------------------------------------------------------------------------
function Split-NumericRun($searchText, $digitPos, $newDigit, $tailStart, $tailEnd, $newTail) {
    $full = $d.Content
    $full.Find.Execute($searchText) | Out-Null
    $start = $full.Start

    $rDigit = $d.Range($start + $digitPos, $start + $digitPos + 1)
    $rDigit.Bold = 1
    $rDigit.Text = $newDigit
    $rDigitBack = $d.Range($start + $digitPos, $start + $digitPos + 1)
    $rDigitBack.Bold = 0

    $rTail = $d.Range($start + $tailStart, $start + $tailEnd)
    $rTail.Bold = 1
    $rTail.Text = $newTail
    $rTailBack = $d.Range($start + $tailStart, $start + $tailStart + $newTail.Length)
    $rTailBack.Bold = 0
}

# ---------------------------------------------------------------------------
# 4) Split run "«6» х  3 = 18" into "«" / "2" / "» х  3 = " / "9"
# ---------------------------------------------------------------------------
Split-NumericRun "«6» х  3 = 18" 1 "2" 11 13 "9"

# ---------------------------------------------------------------------------
# 5) Split run "«8» х 3 = 24" into "«" / "4" / "» х 3 = " / "12"
# ---------------------------------------------------------------------------
Split-NumericRun "«8» х 3 = 24" 1 "4" 10 12 "12"

Write-Output "edits applied"
